$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a daily table of "Betarraga" (beet) prices. Two new rows of
# data need to be inserted right before the current row 298, pushing the
# existing rows 298:382 down to 300:384 (dimension grows from R382 to R384).

$ws.Rows("298:299").Insert()

# Row 298 - new "Primera" quality record
$ws.Range("A298").Value2 = 3
$ws.Range("B298").Value2 = "Femacal de La Calera"
$ws.Range("C298").Value2 = "Coquimbo"
$ws.Range("D298").Value2 = 44463
$ws.Range("E298").Value2 = 5
$ws.Range("F298").Value2 = 100114014
$ws.Range("G298").Value2 = "Betarraga"
$ws.Range("H298").Value2 = "Sin especificar"
$ws.Range("I298").Value2 = "Primera"
$ws.Range("J298").Value2 = 3100
$ws.Range("K298").Value2 = 500
$ws.Range("L298").Value2 = 550
$ws.Range("M298").Value2 = 524
$ws.Range("N298").Value2 = "`$/paquete 4 unidades"
$ws.Range("O298").Value2 = "Provincia de Quillota"
$ws.Range("P298").Value2 = 131
$ws.Range("Q298").Value2 = 4
$ws.Range("R298").Value2 = "Hortaliza"

# Row 299 - new "Segunda" quality record
$ws.Range("A299").Value2 = 3
$ws.Range("B299").Value2 = "Femacal de La Calera"
$ws.Range("C299").Value2 = "Coquimbo"
$ws.Range("D299").Value2 = 44463
$ws.Range("E299").Value2 = 5
$ws.Range("F299").Value2 = 100114014
$ws.Range("G299").Value2 = "Betarraga"
$ws.Range("H299").Value2 = "Sin especificar"
$ws.Range("I299").Value2 = "Segunda"
$ws.Range("J299").Value2 = 1400
$ws.Range("K299").Value2 = 400
$ws.Range("L299").Value2 = 400
$ws.Range("M299").Value2 = 400
$ws.Range("N299").Value2 = "`$/paquete 4 unidades"
$ws.Range("O299").Value2 = "Provincia de Quillota"
$ws.Range("P299").Value2 = 100
$ws.Range("Q299").Value2 = 4
$ws.Range("R299").Value2 = "Hortaliza"

Write-Output "Inserted two rows at 298:299 and populated new data."
